$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = 9.8
$ws.Activate()
$ws.Range("R9").Select()
